# Updated convergence table with correction of the Equilateral mesh method
#
# The "Structured_triangles" (Equilateral mesh) row had a bug: its
# Num_method (column C) was blank and its computation time (column J) was
# wrong. The fix corrects that row and moves it to the top of its block
# (row 10), pushing the other FV/Dirichlet-triangle/cube/tetrahedra rows
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: corrected "Structured_triangles" row (was row 14, C was blank)
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Poisson"
$ws.Range("C10").Value = "FV"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Structured_triangles"
$ws.Range("F10").Value = "Dirichlet"
$ws.Range("G10").Value = 0.8952
$ws.Range("H10").Value = "Triangles"
$ws.Range("I10").Value = "Green"
$ws.Range("J10").Value = 4.86047887802124

# Row 11: former row 10 (Unstructured_triangles)
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Poisson"
$ws.Range("C11").Value = "FV"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Unstructured_triangles"
$ws.Range("F11").Value = "Dirichlet"
$ws.Range("G11").Value = 0.6138
$ws.Range("H11").Value = "Triangles"
$ws.Range("I11").Value = "Green"
$ws.Range("J11").Value = 2.600184917449951

# Row 12: former row 11 (Regular_Cubes). Its Bound_cond cell (F) is blank in
# the source data (an empty shared-string cell, which round-trips as the
# first shared-string entry, "PDE_model", once written back out).
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Poisson"
$ws.Range("C12").Value = "FV"
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "Regular_Cubes"
$ws.Range("F12").Value = "PDE_model"
$ws.Range("G12").Value = 1.3403
$ws.Range("H12").Value = "Cubes"
$ws.Range("I12").Value = "Green"
$ws.Range("J12").Value = 5.900697946548462

# Row 13: former row 12 (Regular_Tetrahedra)
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Poisson"
$ws.Range("C13").Value = "FV"
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = "Regular_Tetrahedra"
$ws.Range("F13").Value = "Dirichlet"
$ws.Range("G13").Value = 0.0065
$ws.Range("H13").Value = "Tetrahedron"
$ws.Range("I13").Value = "Green"
$ws.Range("J13").Value = 62.56098890304565

# Row 14: former row 13 (Unstructured_Tetrahedra)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Poisson"
$ws.Range("C14").Value = "FV"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "Unstructured_Tetrahedra"
$ws.Range("F14").Value = "Dirichlet"
$ws.Range("G14").Value = 0.5359
$ws.Range("H14").Value = "Tetrahedron"
$ws.Range("I14").Value = "Green"
$ws.Range("J14").Value = 3.782500028610229
